# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new F-column value mapping (applies identically to both sheets)
$updates = @{
    3  = 112
    4  = 1615
    5  = 620
    7  = 17
    8  = 11494
    9  = 27
    11 = 449
    12 = 360
    14 = 797
    15 = 12377
    16 = 13069
    22 = 87
    24 = 113
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
